# Insert a new data row above row 134 (pushing the existing rows 134-238
# down to 135-239) and populate the new row 134 with the new record.
#
# Net effect matches the target diff: row 134's old contents now live in
# row 135, row 135's old contents now live in row 136, ... and so on down
# to the old row 238's contents now living in the newly created row 239.
# The freshly inserted row 134 carries the new price-report values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134:238 down to 135:239, inheriting row 134's formatting
# (including the date-time number format on column D) for the new row.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record's data.
$ws.Cells.Item(134, 1).Value = 4
$ws.Cells.Item(134, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(134, 3).Value = "Los Lagos"
$ws.Cells.Item(134, 4).Value = 44673
$ws.Cells.Item(134, 5).Value = 10
$ws.Cells.Item(134, 6).Value = 100112017
$ws.Cells.Item(134, 7).Value = "Apio"
$ws.Cells.Item(134, 8).Value = "Americana (o)"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 40
$ws.Cells.Item(134, 11).Value = 13000
$ws.Cells.Item(134, 12).Value = 14000
$ws.Cells.Item(134, 13).Value = 13500
$ws.Cells.Item(134, 14).Value = "`$/docena de matas"
$ws.Cells.Item(134, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(134, 16).Value = 2250
$ws.Cells.Item(134, 17).Value = 6
$ws.Cells.Item(134, 18).Value = "Hortaliza"
